$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9978650212287903
$ws.Range("B1").Value = 3.43825101852417
$ws.Range("C1").Value = 3.977217435836792
$ws.Range("D1").Value = 3.102188348770142
$ws.Range("E1").Value = 1.314154624938965
